# ImportConflictResolverSequenceDiagram: shift the whole sequence-diagram
# content to the right (the "Model" lifeline column got wider / moved),
# and tighten the loop-guard wording ("from Event" -> "in Event").
#
# All Left/Width values below are expressed in points (PowerPoint COM's
# native unit for Shape geometry; 1 pt = 12700 EMU) and were chosen so
# that, after the host's point->EMU conversion, they reproduce the exact
# EMU offsets/extents from the target OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Rectangle 65 (outer frame) -- moves right AND gets wider (right edge
# shifts by the same amount as everything else, left edge less so).
$sh = $s.Shapes.Item(1)
$sh.Left = 84.0
$sh.Width = 615.646728515625

# Rectangle 35 (lifeline box)
$sh = $s.Shapes.Item(2)
$sh.Left = 204.76110236220472

# Rectangle 62 (":ModelManager" header)
$sh = $s.Shapes.Item(3)
$sh.Left = 184.90867614746094

# Straight Connector 4 (ModelManager lifeline)
$sh = $s.Shapes.Item(4)
$sh.Left = 264.70892333984375

# Rectangle 5 (activation bar)
$sh = $s.Shapes.Item(5)
$sh.Left = 259.03900146484375

# Straight Connector 16 (lifeline)
$sh = $s.Shapes.Item(6)
$sh.Left = 502.72412109375

# Straight Arrow Connector 22 -- also widens like shape 1.
$sh = $s.Shapes.Item(7)
$sh.Left = 48.0
$sh.Width = 211.0389862060547

# TextBox 25 ("handleImportDataAvailableEvent (Event)")
$sh = $s.Shapes.Item(8)
$sh.Left = 84.0

# Straight Arrow Connector 65
$sh = $s.Shapes.Item(9)
$sh.Left = 277.7461417322835

# Rectangle 43 (activation bar)
$sh = $s.Shapes.Item(10)
$sh.Left = 496.3282165527344

# TextBox 56 ("delete(task)")
$sh = $s.Shapes.Item(11)
$sh.Left = 536.37255859375

# Straight Arrow Connector 63
$sh = $s.Shapes.Item(12)
$sh.Left = 280.5482177734375

# Snip Single Corner Rectangle 37 (loop frame corner)
$sh = $s.Shapes.Item(13)
$sh.Left = 205.3002471923828

# TextBox 71 (loop guard text, also edited below)
$sh = $s.Shapes.Item(14)
$sh.Left = 220.0828346456693

# Rectangle 62 (":OverwriteImportConflictResolver" header)
$sh = $s.Shapes.Item(15)
$sh.Left = 394.34905511811024

# Rectangle 40 (activation bar)
$sh = $s.Shapes.Item(16)
$sh.Left = 496.64220472440945

# Rectangle 42 (activation bar)
$sh = $s.Shapes.Item(17)
$sh.Left = 268.92952755905515

# Freeform 44 (arrowhead)
$sh = $s.Shapes.Item(18)
$sh.Left = 273.6476377952756

# Rectangle 46 (activation bar)
$sh = $s.Shapes.Item(19)
$sh.Left = 258.3532283464567

# TextBox 47 ("resolve(add(), delete(), task)")
$sh = $s.Shapes.Item(20)
$sh.Left = 275.86175537109375

# Rectangle 51 (alt frame)
$sh = $s.Shapes.Item(21)
$sh.Left = 211.96087646484375

# Snip Single Corner Rectangle 54 (alt frame corner)
$sh = $s.Shapes.Item(22)
$sh.Left = 211.96102362204724

# TextBox 57 (alt guard text "[task already exists]")
$sh = $s.Shapes.Item(23)
$sh.Left = 223.91867065429688

# Rectangle 48 (activation bar)
$sh = $s.Shapes.Item(24)
$sh.Left = 503.2915954589844

# Freeform 59 (arrowhead)
$sh = $s.Shapes.Item(25)
$sh.Left = 512.3296850393701

# Straight Connector 20 (alt divider line)
$sh = $s.Shapes.Item(26)
$sh.Left = 211.96087646484375

# TextBox 64 ("[else]")
$sh = $s.Shapes.Item(27)
$sh.Left = 238.0828346456693

# TextBox 66 ("add(task)")
$sh = $s.Shapes.Item(28)
$sh.Left = 258.32417322834647

# Rectangle 68 (activation bar)
$sh = $s.Shapes.Item(29)
$sh.Left = 275.19377952755906

# Freeform 69 (arrowhead)
$sh = $s.Shapes.Item(30)
$sh.Left = 284.2318897637795

# TextBox 70 ("add(task)")
$sh = $s.Shapes.Item(31)
$sh.Left = 536.37255859375

# Rectangle 72 (activation bar)
$sh = $s.Shapes.Item(32)
$sh.Left = 503.2915954589844

# Freeform 73 (arrowhead)
$sh = $s.Shapes.Item(33)
$sh.Left = 512.3296850393701

# Wording fix inside the loop guard text box: keep the "loop" run as-is
# and only rewrite the guard-condition run (2nd run) so its formatting
# (the glow effect) is preserved.
$loopShape = $s.Shapes.Item(14)
$loopShape.TextFrame.TextRange.Runs(2).Text = "[for all tasks in Event]"
